$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-18) holds "Förändrad" date values stored as serial 45204.
# Update them to 45205 (advance by one day) to match the source edit.
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
